# This script re-orders the weekly price records (rows 2-17) of the
# "Chirimoya" sheet. Only the columns D (Fecha), L (Calidad), M (Volumen),
# N (Precio mínimo), O (Precio máximo), P (Precio promedio ponderado) and
# S (Precio $/Kg) actually differ between rows, so the edit is expressed
# as: "row N after the edit gets the D/L/M/N/O/P/S values that row M had
# before the edit".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: destination row -> source row (values taken from the ORIGINAL sheet)
$rowMap = @{
    2  = 14
    3  = 15
    4  = 10
    5  = 16
    6  = 17
    7  = 4
    8  = 2
    9  = 12
    10 = 13
    11 = 8
    12 = 9
    13 = 5
    14 = 3
    15 = 6
    16 = 7
    17 = 11
}

# Snapshot the original values for the columns that change, before any
# writes happen, so later writes don't clobber values still needed as a
# source for other rows. NOTE: use `.Value2` (not `.Value`) - this
# interpreter's `.Value` getter does not return the underlying scalar.
$original = @{}
foreach ($r in 2..17) {
    $original[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2   # Fecha
        L = $ws.Cells.Item($r, 12).Value2  # Calidad
        M = $ws.Cells.Item($r, 13).Value2  # Volumen
        N = $ws.Cells.Item($r, 14).Value2  # Precio minimo
        O = $ws.Cells.Item($r, 15).Value2  # Precio maximo
        P = $ws.Cells.Item($r, 16).Value2  # Precio promedio ponderado
        S = $ws.Cells.Item($r, 19).Value2  # Precio $/Kg
    }
}

foreach ($destRow in 2..17) {
    $srcRow = $rowMap[$destRow]
    $vals = $original[$srcRow]

    $ws.Cells.Item($destRow, 4).Value2 = $vals.D
    $ws.Cells.Item($destRow, 12).Value2 = $vals.L
    $ws.Cells.Item($destRow, 13).Value2 = $vals.M
    $ws.Cells.Item($destRow, 14).Value2 = $vals.N
    $ws.Cells.Item($destRow, 15).Value2 = $vals.O
    $ws.Cells.Item($destRow, 16).Value2 = $vals.P
    $ws.Cells.Item($destRow, 19).Value2 = $vals.S
}
